# [Fonds de solidarite] Add 2020-12-31 data
# Update "nombre_aides" (C) and "montant_total" (D) figures for a handful of
# region/classe_effectif rows now that 2020-12-31 data has been folded in.
# The source sheet stores every value as literal text, so each new value is
# written with a leading apostrophe to force Excel to keep it as text
# (matching the existing number-as-text convention) rather than re-typing it
# as a numeric value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($sheet, $addr, $value) {
    $sheet.Range($addr).Value = "'" + $value
}

# Row 35 - Grand Est / 0 salarié
Set-TextValue $ws "C35" "160"
Set-TextValue $ws "D35" "659731.72"

# Row 39 - Grand Est / 10 à 19 salariés
Set-TextValue $ws "C39" "57"
Set-TextValue $ws "D39" "1064717.34"

# Row 65 - Île-de-France / 6 à 9 salariés
Set-TextValue $ws "C65" "1089"
Set-TextValue $ws "D65" "8115492.26"

# Row 90 - Nouvelle-Aquitaine / 0 salarié
Set-TextValue $ws "C90" "289"
Set-TextValue $ws "D90" "1308853.14"

# Row 91 - Nouvelle-Aquitaine / 1 ou 2 salariés
Set-TextValue $ws "C91" "1103"
Set-TextValue $ws "D91" "5860114.39"

# Row 92 - Nouvelle-Aquitaine / 3 à 5 salariés
Set-TextValue $ws "C92" "474"
Set-TextValue $ws "D92" "4393428.38"

# Row 95 - Nouvelle-Aquitaine / 20 à 49 salariés
Set-TextValue $ws "C95" "12"
Set-TextValue $ws "D95" "470000.00"
